# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect that the
# zh-cn and de-de handoffs have been handed back and are in sync with en-US.
# For each locale worksheet it:
#   - Updates the Status column (C) to "Handed back: in sync with en-US"
#   - Adds the "Latest Target File" (F) and "Latest Handback File" (G)
#     columns, mirroring the Source File Name (A) and Latest Handoff File
#     (D) hyperlinked values respectively
#   - Updates the "Latest Handback DateTime" column (H) with the handback
#     timestamp for that locale

$wb = $excel.ActiveWorkbook

function Update-LocaleSheet {
    param(
        [string]$SheetName,
        [string]$HandbackDateTime,
        [string]$Row2MdAddress,
        [string]$Row2XlfAddress,
        [string]$Row3MdAddress,
        [string]$Row3XlfAddress
    )

    $ws = $wb.Worksheets.Item($SheetName)
    Write-Output ("Updating sheet " + $SheetName)

    # --- Status column (C) -------------------------------------------------
    $ws.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
    $ws.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"

    # --- Latest Target File (F) & Latest Handback File (G) -----------------
    $row2MdDisplay = $ws.Cells.Item(2, 1).Text
    $row2XlfDisplay = $ws.Cells.Item(2, 4).Text
    $row3MdDisplay = $ws.Cells.Item(3, 1).Text
    $row3XlfDisplay = $ws.Cells.Item(3, 4).Text

    $ws.Hyperlinks.Add($ws.Cells.Item(2, 6), $Row2MdAddress, "", "", $row2MdDisplay)
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 7), $Row2XlfAddress, "", "", $row2XlfDisplay)
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 6), $Row3MdAddress, "", "", $row3MdDisplay)
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 7), $Row3XlfAddress, "", "", $row3XlfDisplay)

    # --- Latest Handback DateTime (H) --------------------------------------
    $ws.Cells.Item(2, 8).Value = $HandbackDateTime
    $ws.Cells.Item(3, 8).Value = $HandbackDateTime
}

# Overview worksheet: the "zh-cn"/"de-de" status columns (B/C) mirror the
# per-locale Status column, so they need the same text update.
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 2).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(3, 2).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"

# zh-cn worksheet
Update-LocaleSheet "zh-cn" `
    "2016-03-18 10:12:18" `
    "https://github.com/OpenLocalizationTest/oltest/blob/27d485c4c4c3053ee11d33b615e8cf8a6baad9f2/e2e/0d4c1051-cba0-41e3-9001-55bd1a226506.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8deb86d240d9af12f1cc6671b798f5ad493835e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/0d4c1051-cba0-41e3-9001-55bd1a226506.947c82b8c1f32341c9d3d40e48d1ba4f0e56ebc6.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/27d485c4c4c3053ee11d33b615e8cf8a6baad9f2/e2e/893c2691-110f-46fe-897a-8587f67ee692.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8deb86d240d9af12f1cc6671b798f5ad493835e1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/893c2691-110f-46fe-897a-8587f67ee692.5574f8d3747ccc28aa6c0d49cb30783fd2bce155.zh-cn.xlf"

# de-de worksheet
Update-LocaleSheet "de-de" `
    "2016-03-18 10:12:23" `
    "https://github.com/OpenLocalizationTest/oltest/blob/27d485c4c4c3053ee11d33b615e8cf8a6baad9f2/e2e/0d4c1051-cba0-41e3-9001-55bd1a226506.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a34e4af02c0a9db41454e1a44e523d8a16a55bfc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/0d4c1051-cba0-41e3-9001-55bd1a226506.947c82b8c1f32341c9d3d40e48d1ba4f0e56ebc6.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/27d485c4c4c3053ee11d33b615e8cf8a6baad9f2/e2e/893c2691-110f-46fe-897a-8587f67ee692.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a34e4af02c0a9db41454e1a44e523d8a16a55bfc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/893c2691-110f-46fe-897a-8587f67ee692.5574f8d3747ccc28aa6c0d49cb30783fd2bce155.de-de.xlf"
